$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 2 to make room for the "1x6000" entry,
# shifting the existing "4x4090" row down to row 3.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# New row 2: 1x6000
$ws.Cells.Item(2, 1).Value = "1x6000"
$ws.Cells.Item(2, 2).Value = 3168.82
$ws.Cells.Item(2, 3).Value = 1.29
$ws.Cells.Item(2, 4).Value = 0.1130809996570753

# Row 4 (new, after existing 4x4090 row which is now row 3): 4x5090
$ws.Cells.Item(4, 1).Value = "4x5090"
$ws.Cells.Item(4, 2).Value = 4622.08
$ws.Cells.Item(4, 3).Value = 2.6
$ws.Cells.Item(4, 4).Value = 0.1562548078402412
